# TC09_Canine_Filter_Breed-BorderCol.xlsx
# "Fixed variables and query errors in Bread from TC01 to TC30"
#
# The CasesTab (row 2, cell B2) and FilesTab (row 4, cell B4) queries on the
# "startup" sheet had been put in the wrong cells, and the Cases query still
# referenced a `co:cohort` match / `Cohort` output column that doesn't belong
# in this query. This swaps the two query strings into their correct cells
# and drops the stray Cohort match + column from the Cases query.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Border Collie']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Border Collie']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# B2 (CasesTab row) previously held the Files query; it should hold the
# corrected Cases query (no Cohort match/column).
$ws.Range("B2").Value2 = $casesQuery

# B4 (FilesTab row) previously held the (buggy) Cases query; it should hold
# the Files query.
$ws.Range("B4").Value2 = $filesQuery

# The rewrapped query text changes how tall each wrapped/merged row needs to
# be; match the row heights Excel computed after the edit.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Update the saved selection/zoom to match the post-edit view state.
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 115
